$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column appended after the existing "sum" column (G).
# Copy G1's formatting (bold header style with borders/alignment) onto H1
# first, then overwrite its value so the copied style is kept.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Data row: the corresponding "Save" value for row 2.
$ws.Range("H2").Value = 0
